# Generate Report for Handoff
#
# Updates the zh-cn and de-de localization-status sheets: rows 4-7
# (the four files whose Priority was "low") move to "ht", and their
# "Latest Handoff Datetime" timestamps advance by ~30s now that the
# handoff report has (re)run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: Priority (E) low -> ht, Latest Handoff Datetime (H) 04:31:57 -> 04:32:27
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-31 04:32:27"

# de-de: Priority (E) low -> ht, Latest Handoff Datetime (H) 04:32:07 -> 04:32:32
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-31 04:32:32"
